$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.96
$ws.Range("G2").Value = 1.98
$ws.Range("P2").Value = 2.86
$ws.Range("R2").Value = 1.77
$ws.Range("S2").Value = 2.22
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.86
$ws.Range("Y2").Value = 24
$ws.Range("Z2").Value = 34
$ws.Range("AJ2").Value = 24
$ws.Range("AN2").Value = 8
$ws.Range("F3").Value = 1.09
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.01
$ws.Range("J3").Value = 1.09
$ws.Range("V3").Value = 1.13
$ws.Range("W3").Value = 1.01
$ws.Range("G4").Value = 5.5
$ws.Range("I4").Value = 1.73
$ws.Range("J4").Value = 4.4
$ws.Range("K4").Value = 4.9
$ws.Range("F5").Value = 2.18
$ws.Range("I5").Value = 4.8
$ws.Range("K5").Value = 3.45
$ws.Range("L5").Value = 1.53
$ws.Range("AB5").Value = 15
$ws.Range("F6").Value = 1.4
$ws.Range("I6").Value = 14.5
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 4.9
$ws.Range("L6").Value = 1.38
$ws.Range("M6").Value = 1.08
$ws.Range("S6").Value = 4.2
$ws.Range("T6").Value = 2.44
$ws.Range("AB6").Value = 6.8
$ws.Range("R7").Value = 2
$ws.Range("G8").Value = 1.41
$ws.Range("I8").Value = 11
$ws.Range("J8").Value = 5.5
$ws.Range("K8").Value = 7.4
$ws.Range("N8").Value = 6
$ws.Range("O8").Value = 1.16
$ws.Range("Q8").Value = 1.46
$ws.Range("R8").Value = 1.7
$ws.Range("S8").Value = 2.16
$ws.Range("T8").Value = 1.82
$ws.Range("W8").Value = 3.4
$ws.Range("I9").Value = 2.4
$ws.Range("K9").Value = 4
$ws.Range("P9").Value = 1.98
$ws.Range("Q9").Value = 1.83
$ws.Range("U9").Value = 2.18
$ws.Range("V9").Value = 1.71
$ws.Range("J10").Value = 3.7
$ws.Range("P10").Value = 1.92
$ws.Range("Q10").Value = 1.47
$ws.Range("S10").Value = 2.56
$ws.Range("G11").Value = 2.16
$ws.Range("P11").Value = 2.54
$ws.Range("W11").Value = 1.86
$ws.Range("H12").Value = 3.1
$ws.Range("K12").Value = 4.2
$ws.Range("O12").Value = 1.38
$ws.Range("P12").Value = 1.74
$ws.Range("Q12").Value = 1.98
$ws.Range("U12").Value = 1.98
$ws.Range("X12").Value = 1000
$ws.Range("AC12").Value = 1000
$ws.Range("AG12").Value = 1000
$ws.Range("H13").Value = 1.45
$ws.Range("K13").Value = 5.6
$ws.Range("N13").Value = 5.3
$ws.Range("O13").Value = 1.19
$ws.Range("Q13").Value = 1.56
$ws.Range("T13").Value = 1.76
$ws.Range("U13").Value = 2.08
$ws.Range("H14").Value = 3.2
$ws.Range("L14").Value = 1.29
$ws.Range("P14").Value = 1.98
$ws.Range("Q14").Value = 1.81
$ws.Range("W14").Value = 1.73
$ws.Range("X14").Value = 990
$ws.Range("Y14").Value = 17.5
$ws.Range("AB14").Value = 13.5
$ws.Range("AC14").Value = 10.5
$ws.Range("AD14").Value = 17
$ws.Range("AG14").Value = 13
$ws.Range("H15").Value = 10.5
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 7.4
$ws.Range("K15").Value = 9.199999999999999
$ws.Range("P15").Value = 3.2
$ws.Range("R15").Value = 1.88
$ws.Range("T15").Value = 1.81
$ws.Range("U15").Value = 2
$ws.Range("AC15").Value = 22
$ws.Range("AB16").Value = 12
$ws.Range("F17").Value = 5.3
$ws.Range("G17").Value = 6.2
$ws.Range("H17").Value = 1.64
$ws.Range("I17").Value = 1.73
$ws.Range("J17").Value = 3.95
$ws.Range("K17").Value = 4.5
$ws.Range("H18").Value = 2.82
$ws.Range("J18").Value = 3.2
$ws.Range("K18").Value = 3.55
$ws.Range("L18").Value = 1.42
$ws.Range("Q18").Value = 2.04
$ws.Range("W18").Value = 1.56
$ws.Range("AK18").Value = 980
$ws.Range("AL18").Value = 980
$ws.Range("P19").Value = 2.24
$ws.Range("W19").Value = 1.2
$ws.Range("X19").Value = 28
$ws.Range("F20").Value = 5.8
$ws.Range("G20").Value = 7.6
$ws.Range("H20").Value = 1.47
$ws.Range("I20").Value = 1.55
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 6
$ws.Range("L20").Value = 1.16
$ws.Range("M20").Value = 1.01
$ws.Range("N20").Value = 6.8
$ws.Range("P20").Value = 2.96
$ws.Range("R20").Value = 1.81
$ws.Range("T20").Value = 1.49
$ws.Range("V20").Value = 2.72
$ws.Range("X20").Value = 990
$ws.Range("Y20").Value = 990
$ws.Range("AC20").Value = 990
$ws.Range("AF20").Value = 65
$ws.Range("AJ20").Value = 180
$ws.Range("AK20").Value = 70
$ws.Range("AL20").Value = 70
$ws.Range("AM20").Value = 70
$ws.Range("AN20").Value = 60
$ws.Range("AO20").Value = 5
$ws.Range("G21").Value = 5.8
$ws.Range("H21").Value = 1.73
$ws.Range("N21").Value = 3.7
$ws.Range("O21").Value = 1.3
$ws.Range("W21").Value = 1.22
$ws.Range("X21").Value = 18.5
$ws.Range("Q22").Value = 2.46
$ws.Range("F23").Value = 2.16
$ws.Range("G23").Value = 2.44
$ws.Range("H23").Value = 2.94
$ws.Range("I23").Value = 3.45
$ws.Range("J23").Value = 3.75
$ws.Range("K23").Value = 4.5
$ws.Range("L23").Value = 1.2
$ws.Range("M23").Value = 1.03
$ws.Range("N23").Value = 5.1
$ws.Range("O23").Value = 1.19
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.42
$ws.Range("U23").Value = 2.5
$ws.Range("V23").Value = 1.41
$ws.Range("W23").Value = 1.7
$ws.Range("X23").Value = 29
$ws.Range("Y23").Value = 22
$ws.Range("Z23").Value = 32
$ws.Range("AB23").Value = 17.5
$ws.Range("AC23").Value = 10.5
$ws.Range("AD23").Value = 17.5
$ws.Range("AE23").Value = 38
$ws.Range("AF23").Value = 22
$ws.Range("AG23").Value = 12.5
$ws.Range("AH23").Value = 18.5
$ws.Range("AI23").Value = 42
$ws.Range("AJ23").Value = 32
$ws.Range("AK23").Value = 23
$ws.Range("AL23").Value = 36
$ws.Range("AN23").Value = 14.5
$ws.Range("AO23").Value = 22
$ws.Range("F24").Value = 1.21
$ws.Range("G24").Value = 1.23
$ws.Range("J24").Value = 8
$ws.Range("N24").Value = 10
$ws.Range("P24").Value = 4
$ws.Range("Q24").Value = 1.26
$ws.Range("R24").Value = 2.28
$ws.Range("S24").Value = 1.66
$ws.Range("T24").Value = 1.65
$ws.Range("U24").Value = 2.24
$ws.Range("W24").Value = 5
$ws.Range("X24").Value = 75
$ws.Range("Z24").Value = 220
$ws.Range("AA24").Value = 630
$ws.Range("AB24").Value = 22
$ws.Range("AE24").Value = 200
$ws.Range("AF24").Value = 13
$ws.Range("AG24").Value = 15.5
$ws.Range("AI24").Value = 140
$ws.Range("AJ24").Value = 14.5
$ws.Range("AK24").Value = 13
$ws.Range("AL24").Value = 28
$ws.Range("AM24").Value = 130
$ws.Range("AN24").Value = 2.8
$ws.Range("N25").Value = 5.2
$ws.Range("P25").Value = 2.44
$ws.Range("Q25").Value = 1.68
$ws.Range("U25").Value = 2.74
$ws.Range("AI25").Value = 30
$ws.Range("AK25").Value = 23
$ws.Range("AL25").Value = 29
$ws.Range("AM25").Value = 50
$ws.Range("AO25").Value = 18
